# Append a new auto-logged experiment row (row 39) to the experiment log,
# mirroring the existing "Auto-log: Q=6, D=4, Skip=add" row but trained for
# more epochs (best model so far).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

$ws.Cells.Item($row, 1).Value2  = 38                                     # A: Experiment ID
$ws.Cells.Item($row, 2).Value2  = "Auto-log: Q=6, D=4, Skip=add"         # B: Description
$ws.Cells.Item($row, 3).Value2  = 0                                      # C: LSTM Layers
$ws.Cells.Item($row, 4).Value2  = 4                                      # D: GRU Layers
$ws.Cells.Item($row, 5).Value2  = 0                                      # E: RNN Layers
$ws.Cells.Item($row, 6).Value2  = 32                                     # F: Hidden Size
$ws.Cells.Item($row, 7).Value2  = 20                                     # G: Window Size
$ws.Cells.Item($row, 8).Value2  = $true                                  # H: Use Quantum
$ws.Cells.Item($row, 9).Value2  = 6                                      # I: Qubits
$ws.Cells.Item($row, 10).Value2 = 4                                      # J: Q Depth
$ws.Cells.Item($row, 11).Value2 = 3                                      # K: Rotation Params
$ws.Cells.Item($row, 12).Value2 = "add"                                  # L: Skip Connection
$ws.Cells.Item($row, 13).Value2 = "Tanh"                                 # M: Post-Quantum Activation
$ws.Cells.Item($row, 14).Value2 = $false                                 # N: Use Dropout
$ws.Cells.Item($row, 15).Value2 = 0                                      # O: Dropout Rate
$ws.Cells.Item($row, 16).Value2 = $false                                 # P: Use LayerNorm
$ws.Cells.Item($row, 17).Value2 = ""                                     # Q: Final Activation
$ws.Cells.Item($row, 18).Value2 = 40                                     # R: Num Epochs
$ws.Cells.Item($row, 19).Value2 = 10                                     # S: Early Stop Patience
$ws.Cells.Item($row, 20).Value2 = 0.00008380591476994363                 # T: Train Loss
$ws.Cells.Item($row, 21).Value2 = 0.002030752476012493                   # U: Validation Loss
$ws.Cells.Item($row, 22).Value2 = 2.565991401672363                      # V: MAE
$ws.Cells.Item($row, 23).Value2 = 3.171932220458984                      # W: RMSE
$ws.Cells.Item($row, 24).Value2 = 1.774586319923401                      # X: Avg % Error
$ws.Cells.Item($row, 25).Value2 = "test run with classical only layer"   # Y: Notes
$ws.Cells.Item($row, 26).Value2 = "AAPL, MSFT, GOOGL"                    # Z: Tickers
